# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets to the newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1864
$ws1.Range("G2").Value = 60
$ws1.Range("F3").Value = 488
$ws1.Range("F6").Value = 2550
$ws1.Range("F7").Value = 168
$ws1.Range("F10").Value = 1523
$ws1.Range("F11").Value = 527
$ws1.Range("F12").Value = 42
$ws1.Range("F13").Value = 327
$ws1.Range("F14").Value = 229
$ws1.Range("F17").Value = 207
$ws1.Range("F20").Value = 12
$ws1.Range("F21").Value = 171
$ws1.Range("F22").Value = 55
$ws1.Range("F23").Value = 1623
$ws1.Range("F24").Value = 25
$ws1.Range("F25").Value = 397
$ws1.Range("F26").Value = 571
$ws1.Range("F28").Value = 296
$ws1.Range("F29").Value = 412

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1864
$ws4.Range("G2").Value = 60
$ws4.Range("F4").Value = 488
$ws4.Range("F7").Value = 2550
$ws4.Range("F8").Value = 168
$ws4.Range("F11").Value = 1523
$ws4.Range("F12").Value = 527
$ws4.Range("F13").Value = 42
$ws4.Range("F14").Value = 327
$ws4.Range("F15").Value = 229
$ws4.Range("F18").Value = 207
$ws4.Range("F21").Value = 12
$ws4.Range("F22").Value = 171
$ws4.Range("F23").Value = 55
$ws4.Range("F24").Value = 1623
$ws4.Range("F25").Value = 25
$ws4.Range("F26").Value = 397
$ws4.Range("F27").Value = 571
$ws4.Range("F29").Value = 296
$ws4.Range("F30").Value = 412
